$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.015.98"

# Row 3
$ws.Range("D3").Value = "3.657.12"
$ws.Range("E3").Value = "  +5.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "595.61"
$ws.Range("E5").Value = "  +1.91%  "

# Row 6
$ws.Range("D6").Value = "195.43"
$ws.Range("E6").Value = "  +4.97%  "

# Row 7
$ws.Range("D7").Value = "0.649"
$ws.Range("E7").Value = "  +2.76%  "

# Row 8
$ws.Range("D8").Value = "3.649.25"
$ws.Range("E8").Value = "  +5.83%  "

# Row 9
$ws.Range("E9").Value = "  +0.14%  "

# Row 10
$ws.Range("D10").Value = "0.181"

# Row 11
$ws.Range("D11").Value = "0.675"
$ws.Range("E11").Value = "  +4.98%  "

# Row 12
$ws.Range("D12").Value = "58.76"
$ws.Range("E12").Value = "  +4.69%  "

# Row 13
$ws.Range("D13").Value = "0.0000294"
$ws.Range("E13").Value = "  +6.21%  "

# Row 14
$ws.Range("D14").Value = "9.98"
$ws.Range("E14").Value = "  +6.55%  "

# Row 15
$ws.Range("D15").Value = "4.251.48"
$ws.Range("E15").Value = "  +5.76%  "

# Row 16
$ws.Range("D16").Value = "19.94"
$ws.Range("E16").Value = "  +6.81%  "

# Row 17
$ws.Range("D17").Value = "3.659.66"
$ws.Range("E17").Value = "  +5.66%  "

# Row 18
$ws.Range("D18").Value = "71.060.67"
$ws.Range("E18").Value = "  +6.10%  "

# Row 19
$ws.Range("D19").Value = "12.82"
$ws.Range("E19").Value = "  +6.22%  "

# Row 20
$ws.Range("D20").Value = "0.121"
$ws.Range("E20").Value = "  +3.28%  "

# Row 21
$ws.Range("E21").Value = "  +5.16%  "

# Row 22
$ws.Range("D22").Value = "492.14"
$ws.Range("E22").Value = "  +0.38%  "

# Row 23
$ws.Range("D23").Value = "19.07"
$ws.Range("E23").Value = "  +15.43%  "

# Row 24
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").Value = "4.52"
$ws.Range("E25").Value = "  +2.28%  "

# Row 26
$ws.Range("D26").Value = "91.54"
$ws.Range("E26").Value = "  +2.28%  "

# Row 27
$ws.Range("D27").Value = "3.17"
$ws.Range("E27").Value = "  +8.41%  "

# Row 28
$ws.Range("D28").Value = "11.61"
$ws.Range("E28").Value = "  +7.03%  "

# Row 29
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  +7.12%  "

# Row 30
$ws.Range("D30").Value = "32.94"
$ws.Range("E30").Value = "  +5.56%  "

# Row 31
$ws.Range("E31").Value = "  +11.71%  "

# Row 32
$ws.Range("E32").Value = "  +9.65%  "

# Row 33
$ws.Range("D33").Value = "633.34"
$ws.Range("E33").Value = "  +5.60%  "

# Row 34
$ws.Range("D34").Value = "12.30"
$ws.Range("E34").Value = "  +5.42%  "

# Row 35
$ws.Range("D35").Value = "65.72"
$ws.Range("E35").Value = "  +3.64%  "

# Row 36
$ws.Range("D36").Value = "40.88"
$ws.Range("E36").Value = "  +12.53%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0835"
$ws.Range("E37").Value = "  +11.28%  "

# Row 38
$ws.Range("D38").Value = "0.415"
$ws.Range("E38").Value = "  +8.53%  "

# Row 39
$ws.Range("E39").Value = "  -0.90%  "

# Row 40
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.07%  "

# Row 41
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  +2.38%  "

# Row 42
$ws.Range("D42").Value = "3.325.01"
$ws.Range("E42").Value = "  +2.83%  "

# Row 43
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.88"
$ws.Range("E43").Value = "  +15.30%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "3.17"
$ws.Range("E44").Value = "  +9.79%  "

# Row 45
$ws.Range("E45").Value = "  +6.81%  "

# Row 46
$ws.Range("E46").Value = "  +4.79%  "

# Row 47
$ws.Range("D47").Value = "0.139"
$ws.Range("E47").Value = "  +3.59%  "

# Row 48
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "9.32"
$ws.Range("E48").Value = "  +7.22%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.29"
$ws.Range("E49").Value = "  +2.00%  "

# Row 50
$ws.Range("D50").Value = "3.33"
$ws.Range("E50").Value = "  +1.87%  "

# Row 51
$ws.Range("E51").Value = "  -0.17%  "
